$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $found = $range.Find.Execute(
        $findText,    # FindText
        $false,       # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replaceText, # ReplaceWith
        2             # Replace (wdReplaceOne)
    )
    if (-not $found) {
        Write-Output "WARNING: text not found: $findText"
    }
}

# --- "Programa resumido" paragraph: split into two lines with a manual line break ---
Replace-Text `
    "Introdução a segurança e medicina do trabalho; Conceitos técnico e aspectos legais em biossegurança;" `
    "Introdução a segurança e medicina do trabalho; ^lConceitos técnico e aspectos legais em biossegurança;"

# --- "Programa" paragraph: split the run-on list into one topic per line ---
Replace-Text `
    "Introdução a gestão em segurança do trabalho e estratégias de prevenção; Identificação e controle dos riscos ambientais (físicos, químicos e biológicos). Normas regulamentadoras;Classes de risco biológico, níveis de biossegurança e normas para a atividades de biotecnologia. Descarte e classificação de resíduo;Legislação para produção e manejo organismos geneticamente modificados (OGM) e seus derivados;Biossegurança no manuseio de cobaias; Princípios de bioética;Estudos de casos problemas e soluções" `
    "Introdução a gestão em segurança do trabalho e estratégias de prevenção; Identificação e controle dos riscos ambientais (físicos, químicos e biológicos). Normas regulamentadoras;^lClasses de risco biológico, níveis de biossegurança e normas para a atividades de biotecnologia. ^lDescarte e classificação de resíduo;^lLegislação para produção e manejo organismos geneticamente modificados (OGM) e seus derivados;^lBiossegurança no manuseio de cobaias; ^lPrincípios de bioética;^lEstudos de casos problemas e soluções"

# --- "Bibliografia" paragraph: split the four references, one per (double-spaced) line ---
Replace-Text `
    "1-Binsfeld, P. C. Fundamentos Técnicos e o Sistema Nacional de Biossegurança em Biotecnologia. Interciência, 1ª edição 2015.2-Gonçalves Simão, L. B. Gestão de Segurança e Medicina do Trabalho, Normas Regulamentadoras e Fator Acidentário de Prevenção. Cenofisco, 1ª edição 2015.3-Hirata, M.H., Mancini Filho, J. Hirata, R. D. C. Manual de biossegurança.  Editora Manole. 3ª edição 2016.4- Semplici, S. Onze Teses de Bioética. Editora Ideias e Letras;1ª edição 2014" `
    "1-Binsfeld, P. C. Fundamentos Técnicos e o Sistema Nacional de Biossegurança em Biotecnologia. Interciência, 1ª edição 2015.^l^l2-Gonçalves Simão, L. B. Gestão de Segurança e Medicina do Trabalho, Normas Regulamentadoras e Fator Acidentário de Prevenção. Cenofisco, 1ª edição 2015.^l^l3-Hirata, M.H., Mancini Filho, J. Hirata, R. D. C. Manual de biossegurança.  Editora Manole. 3ª edição 2016.^l^l4- Semplici, S. Onze Teses de Bioética. Editora Ideias e Letras;1ª edição 2014"
